$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 894.3077
$ws.Range("I41").Value = 666.8333
$ws.Range("J41").Value = 1089.2858
$ws.Range("K41").Value = 666.8333
$ws.Range("L41").Value = 1089.2858
$ws.Range("M41").Value = -226.8333
$ws.Range("N41").Value = -1969.2858
$ws.Range("H58").Value = 29279.182
$ws.Range("I58").Value = 31207.1
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 93621.29999999999
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -93471.29999999999
$ws.Range("N58").Value = -30300
$ws.Range("H62").Value = 3059.9375
$ws.Range("I62").Value = 3104.2856
$ws.Range("J62").Value = 2749.5
$ws.Range("K62").Value = 3104.2856
$ws.Range("L62").Value = 2749.5
$ws.Range("M62").Value = -2480.2856
$ws.Range("N62").Value = -3997.5
$ws.Range("H64").Value = 3454.6365
$ws.Range("I64").Value = 3166.7778
$ws.Range("J64").Value = 4750
$ws.Range("K64").Value = 3166.7778
$ws.Range("L64").Value = 4750
$ws.Range("M64").Value = -2918.7778
$ws.Range("N64").Value = -5246
$ws.Range("H65").Value = 3059.9375
$ws.Range("I65").Value = 3104.2856
$ws.Range("J65").Value = 2749.5
$ws.Range("K65").Value = 15521.428
$ws.Range("L65").Value = 13747.5
$ws.Range("M65").Value = -12401.428
$ws.Range("N65").Value = -19987.5
$ws.Range("H67").Value = 3454.6365
$ws.Range("I67").Value = 3166.7778
$ws.Range("J67").Value = 4750
$ws.Range("K67").Value = 3166.7778
$ws.Range("L67").Value = 4750
$ws.Range("M67").Value = -2308.7778
$ws.Range("N67").Value = -6466
$ws.Range("H97").Value = 200171
$ws.Range("J97").Value = 200171
$ws.Range("L97").Value = 600513
$ws.Range("N97").Value = -601505
$ws.Range("H98").Value = 1495.7059
$ws.Range("I98").Value = 1514.1875
$ws.Range("K98").Value = 1514.1875
$ws.Range("M98").Value = -16.1875
$ws.Range("H103").Value = 5810
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 4255.1724
$ws.Range("I106").Value = 3420
$ws.Range("J106").Value = 6111.1113
$ws.Range("K106").Value = 3420
$ws.Range("L106").Value = 6111.1113
$ws.Range("M106").Value = -2789
$ws.Range("N106").Value = -7373.1113
$ws.Range("H109").Value = 30500
$ws.Range("J109").Value = 30500
$ws.Range("L109").Value = 30500
$ws.Range("N109").Value = -33274
$ws.Range("H112").Value = 1261.9166
$ws.Range("J112").Value = 1359.6
$ws.Range("L112").Value = 4078.8
$ws.Range("N112").Value = -6294.799999999999
$ws.Range("H115").Value = 854.75
$ws.Range("I115").Value = 854.75
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2564.25
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -997.25
$ws.Range("H118").Value = 689.3333
$ws.Range("I118").Value = 334.2857
$ws.Range("K118").Value = 1002.8571
$ws.Range("M118").Value = 654.1428999999999
$ws.Range("H122").Value = 1495.7059
$ws.Range("I122").Value = 1514.1875
$ws.Range("K122").Value = 4542.5625
$ws.Range("M122").Value = -2092.5625
$ws.Range("H127").Value = 50001064
$ws.Range("I127").Value = 166667260
$ws.Range("J127").Value = 1265.3572
$ws.Range("K127").Value = 500001780
$ws.Range("L127").Value = 3796.0716
$ws.Range("M127").Value = -499996820
$ws.Range("N127").Value = -13716.0716
$ws.Range("H129").Value = 911.3125
$ws.Range("I129").Value = 233.4
$ws.Range("J129").Value = 1036.8518
$ws.Range("K129").Value = 700.2
$ws.Range("L129").Value = 3110.5554
$ws.Range("M129").Value = 4299.8
$ws.Range("N129").Value = -13110.5554
$ws.Range("H131").Value = 142859870
$ws.Range("I131").Value = 200001070
$ws.Range("J131").Value = 6900
$ws.Range("K131").Value = 600003210
$ws.Range("L131").Value = 20700
$ws.Range("M131").Value = -599998170
$ws.Range("N131").Value = -30780
$ws.Range("H135").Value = 827.8889
$ws.Range("I135").Value = 801.875
$ws.Range("J135").Value = 1036
$ws.Range("K135").Value = 7216.875
$ws.Range("L135").Value = 9324
$ws.Range("M135").Value = -4681.875
$ws.Range("N135").Value = -14394
$ws.Range("H137").Value = 655344.3
$ws.Range("I137").Value = 1506.9697
$ws.Range("K137").Value = 4520.909100000001
$ws.Range("M137").Value = -1970.909100000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1069.8422
$ws.Range("I61").Value = 897.3125
$ws.Range("J61").Value = 1990
$ws.Range("K61").Value = 897.3125
$ws.Range("L61").Value = 1990
$ws.Range("M61").Value = -685.3125
$ws.Range("N61").Value = -2414
$ws.Range("H74").Value = 51014.2
$ws.Range("I74").Value = 59710.824
$ws.Range("J74").Value = 1733.3334
$ws.Range("K74").Value = 59710.824
$ws.Range("L74").Value = 1733.3334
$ws.Range("M74").Value = -58836.824
$ws.Range("N74").Value = -3481.3334
$ws.Range("H77").Value = 51014.2
$ws.Range("I77").Value = 59710.824
$ws.Range("J77").Value = 1733.3334
$ws.Range("K77").Value = 298554.12
$ws.Range("L77").Value = 8666.666999999999
$ws.Range("M77").Value = -294186.12
$ws.Range("N77").Value = -17402.667
$ws.Range("H136").Value = 1069.8422
$ws.Range("I136").Value = 897.3125
$ws.Range("J136").Value = 1990
$ws.Range("K136").Value = 2691.9375
$ws.Range("L136").Value = 5970
$ws.Range("M136").Value = -141.9375
$ws.Range("N136").Value = -11070

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 454
$ws.Range("I11").Value = 144.8
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 144.8
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = -4.800000000000011
$ws.Range("N11").Value = -2280
$ws.Range("H86").Value = 1726.6
$ws.Range("I86").Value = 1533.25
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1533.25
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -410.25
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 1726.6
$ws.Range("I89").Value = 1533.25
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 7666.25
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -2050.25
$ws.Range("N89").Value = -23732
$ws.Range("H134").Value = 4322.768
$ws.Range("I134").Value = 4525.6216
$ws.Range("J134").Value = 3927.7368
$ws.Range("K134").Value = 13576.8648
$ws.Range("L134").Value = 11783.2104
$ws.Range("M134").Value = -11041.8648
$ws.Range("N134").Value = -16853.2104

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13515062
$ws.Range("I31").Value = 1322.5714
$ws.Range("K31").Value = 1322.5714
$ws.Range("M31").Value = -1027.5714
$ws.Range("H34").Value = 13515062
$ws.Range("I34").Value = 1322.5714
$ws.Range("K34").Value = 1322.5714
$ws.Range("M34").Value = -1120.5714

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 25552.5
$ws.Range("J13").Value = 49800
$ws.Range("L13").Value = 49800
$ws.Range("N13").Value = -50080
$ws.Range("H136").Value = 4897.125
$ws.Range("I136").Value = 10650.8
$ws.Range("J136").Value = 2281.818
$ws.Range("K136").Value = 31952.4
$ws.Range("L136").Value = 6845.454000000001
$ws.Range("M136").Value = -29402.4
$ws.Range("N136").Value = -11945.454
